$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-06 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-07 Thursday", 2)
$d.Content.Find.Execute("58×39=2262", $true, $false, $false, $false, $false, $true, 1, $false, "46×15=690", 2)
$d.Content.Find.Execute("50×57=2850", $true, $false, $false, $false, $false, $true, 1, $false, "62×24=1488", 2)
$d.Content.Find.Execute("16×49=784", $true, $false, $false, $false, $false, $true, 1, $false, "83×67=5561", 2)
$d.Content.Find.Execute("14×53=742", $true, $false, $false, $false, $false, $true, 1, $false, "40×93=3720", 2)
$d.Content.Find.Execute("12×60=720", $true, $false, $false, $false, $false, $true, 1, $false, "95×48=4560", 2)
$d.Content.Find.Execute("44×33=1452", $true, $false, $false, $false, $false, $true, 1, $false, "19×25=475", 2)
$d.Content.Find.Execute("39×81=3159", $true, $false, $false, $false, $false, $true, 1, $false, "61×84=5124", 2)
$d.Content.Find.Execute("84×99=8316", $true, $false, $false, $false, $false, $true, 1, $false, "19×92=1748", 2)
$d.Content.Find.Execute("91×24=2184", $true, $false, $false, $false, $false, $true, 1, $false, "21×89=1869", 2)
$d.Content.Find.Execute("31×84=2604", $true, $false, $false, $false, $false, $true, 1, $false, "74×33=2442", 2)
$d.Content.Find.Execute("37×79=2923", $true, $false, $false, $false, $false, $true, 1, $false, "17×84=1428", 2)
$d.Content.Find.Execute("13×19=247", $true, $false, $false, $false, $false, $true, 1, $false, "72×62=4464", 2)
$d.Content.Find.Execute("57×80=4560", $true, $false, $false, $false, $false, $true, 1, $false, "68×20=1360", 2)
$d.Content.Find.Execute("92×29=2668", $true, $false, $false, $false, $false, $true, 1, $false, "31×52=1612", 2)
$d.Content.Find.Execute("61×63=3843", $true, $false, $false, $false, $false, $true, 1, $false, "33×78=2574", 2)
$d.Content.Find.Execute("68×35=2380", $true, $false, $false, $false, $false, $true, 1, $false, "50×29=1450", 2)
$d.Content.Find.Execute("70×32=2240", $true, $false, $false, $false, $false, $true, 1, $false, "29×42=1218", 2)
$d.Content.Find.Execute("76×97=7372", $true, $false, $false, $false, $false, $true, 1, $false, "32×59=1888", 2)
$d.Content.Find.Execute("88×41=3608", $true, $false, $false, $false, $false, $true, 1, $false, "40×21=840", 2)
$d.Content.Find.Execute("65×67=4355", $true, $false, $false, $false, $false, $true, 1, $false, "53×14=742", 2)
$d.Content.Find.Execute("91×25=2275", $true, $false, $false, $false, $false, $true, 1, $false, "92×90=8280", 2)
$d.Content.Find.Execute("23×42=966", $true, $false, $false, $false, $false, $true, 1, $false, "80×75=6000", 2)
$d.Content.Find.Execute("76×88=6688", $true, $false, $false, $false, $false, $true, 1, $false, "34×28=952", 2)
$d.Content.Find.Execute("56×20=1120", $true, $false, $false, $false, $false, $true, 1, $false, "69×14=966", 2)
$d.Content.Find.Execute("41×96=3936", $true, $false, $false, $false, $false, $true, 1, $false, "37×77=2849", 2)
